# Update CDA Logical model for ST.r2b
#
# The workbook has two sheets:
#   1. "Metadata" - a Property/Value table describing the StructureDefinition
#   2. "Elements" - the element-level detail grid
#
# This edit:
#   - bumps the Version metadata value
#   - bumps the Date metadata value
#   - inserts a new "Jurisdiction" property row (empty value) right after
#     "Contact" and before "Description", shifting the rows below it down

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Version and Date values on the Metadata sheet ---
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# --- Insert a new "Jurisdiction" row after "Contact" (row 10) ---
$ws1.Rows.Item(11).Insert()

# Reuse the formatting from the row above (Contact) so the new row keeps
# the same style as every other data row in the table.
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$ws1.Application.CutCopyMode = 0

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

Write-Output "Metadata sheet updated"
